$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $origStyle = $cell.Style
    $cell.Value = "'" + $text
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") '29.876.64'
Set-TextValue $ws.Range("E2") '  -1.86%  '
Set-TextValue $ws.Range("D3") '1.887.22'
Set-TextValue $ws.Range("E3") '  -2.79%  '
Set-TextValue $ws.Range("D4") '1.001'
Set-TextValue $ws.Range("E4") '  -0.48%  '
Set-TextValue $ws.Range("D5") '0.7330'
Set-TextValue $ws.Range("E5") '  -2.90%  '
Set-TextValue $ws.Range("D6") '242.36'
Set-TextValue $ws.Range("E6") '  -1.74%  '
Set-TextValue $ws.Range("D7") '1.001'
Set-TextValue $ws.Range("E7") '  -0.54%  '
Set-TextValue $ws.Range("D8") '0.3101'
Set-TextValue $ws.Range("E8") '  -2.89%  '
Set-TextValue $ws.Range("D9") '26.16'
Set-TextValue $ws.Range("E9") '  -5.52%  '
Set-TextValue $ws.Range("D10") '0.06892'
Set-TextValue $ws.Range("E10") '  -1.58%  '
Set-TextValue $ws.Range("D11") '0.7706'
Set-TextValue $ws.Range("E11") '  -1.45%  '
Set-TextValue $ws.Range("D12") '0.07939'
Set-TextValue $ws.Range("E12") '  -1.04%  '
Set-TextValue $ws.Range("D13") '1.875.13'
Set-TextValue $ws.Range("E13") '  -3.37%  '
Set-TextValue $ws.Range("D14") '5.216'
Set-TextValue $ws.Range("E14") '  -2.56%  '
Set-TextValue $ws.Range("D15") '91.20'
Set-TextValue $ws.Range("E15") '  -3.51%  '
Set-TextValue $ws.Range("D16") '14.17'
Set-TextValue $ws.Range("E16") '  -1.95%  '
Set-TextValue $ws.Range("D17") '29.907.05'
Set-TextValue $ws.Range("E17") '  -1.77%  '
Set-TextValue $ws.Range("D18") '5.742'
Set-TextValue $ws.Range("E18") '  -0.54%  '
Set-TextValue $ws.Range("D19") '239.17'
Set-TextValue $ws.Range("E19") '  -6.01%  '
Set-TextValue $ws.Range("D20") '0.000007748'
Set-TextValue $ws.Range("E20") '  -2.28%  '
Set-TextValue $ws.Range("D21") '1.000'
Set-TextValue $ws.Range("E21") '  -0.39%  '
Set-TextValue $ws.Range("D22") '2.113.69'
Set-TextValue $ws.Range("E22") '  -3.46%  '
Set-TextValue $ws.Range("D23") '1.001'
Set-TextValue $ws.Range("E23") '  -0.52%  '
Set-TextValue $ws.Range("D24") '6.895'
Set-TextValue $ws.Range("E24") '  +3.24%  '
Set-TextValue $ws.Range("D25") '9.287'
Set-TextValue $ws.Range("E25") '  -2.53%  '
Set-TextValue $ws.Range("D26") '164.36'
Set-TextValue $ws.Range("E26") '  -0.60%  '
Set-TextValue $ws.Range("D27") '18.83'
Set-TextValue $ws.Range("E27") '  -1.15%  '
Set-TextValue $ws.Range("D28") '0.1267'
Set-TextValue $ws.Range("E28") '  -4.73%  '
Set-TextValue $ws.Range("D29") '2.010'
Set-TextValue $ws.Range("E29") '  -11.33%  '
Set-TextValue $ws.Range("D30") '1.352'
Set-TextValue $ws.Range("E30") '  -1.80%  '
Set-TextValue $ws.Range("D31") '1.531'
Set-TextValue $ws.Range("E31") '  +1.03%  '
Set-TextValue $ws.Range("D32") '4.297'
Set-TextValue $ws.Range("E32") '  -2.11%  '
Set-TextValue $ws.Range("D33") '4.075'
Set-TextValue $ws.Range("E33") '  -1.09%  '
Set-TextValue $ws.Range("D34") '0.05094'
Set-TextValue $ws.Range("E34") '  -1.25%  '
Set-TextValue $ws.Range("D35") '1.274'
Set-TextValue $ws.Range("E35") '  -0.35%  '
Set-TextValue $ws.Range("D36") '0.7350'
Set-TextValue $ws.Range("E36") '  -1.39%  '
Set-TextValue $ws.Range("D37") '2.722'
Set-TextValue $ws.Range("E37") '  -2.67%  '
Set-TextValue $ws.Range("D38") '0.01916'
Set-TextValue $ws.Range("E38") '  -1.54%  '
Set-TextValue $ws.Range("D39") '2.774'
Set-TextValue $ws.Range("E39") '  -1.35%  '
Set-TextValue $ws.Range("D40") '6.292'
Set-TextValue $ws.Range("E40") '  -1.95%  '
Set-TextValue $ws.Range("D41") '74.12'
Set-TextValue $ws.Range("E41") '  -6.06%  '
Set-TextValue $ws.Range("D42") '0.4450'
Set-TextValue $ws.Range("E42") '  -0.80%  '
Set-TextValue $ws.Range("D43") '1.929'
Set-TextValue $ws.Range("E43") '  -1.90%  '
Set-TextValue $ws.Range("D44") '1.001'
Set-TextValue $ws.Range("E44") '  -0.57%  '
Set-TextValue $ws.Range("D45") '0.8373'
Set-TextValue $ws.Range("E45") '  +0.50%  '
Set-TextValue $ws.Range("D46") '7.623'
Set-TextValue $ws.Range("E46") '  +1.95%  '
Set-TextValue $ws.Range("D47") '100.70'
Set-TextValue $ws.Range("E47") '  -0.48%  '
Set-TextValue $ws.Range("D48") '9.783'
Set-TextValue $ws.Range("E48") '  +0.35%  '
Set-TextValue $ws.Range("B49") 'Elrond'
Set-TextValue $ws.Range("C49") 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextValue $ws.Range("D49") '36.90'
Set-TextValue $ws.Range("E49") '  -0.78%  '
Set-TextValue $ws.Range("B50") 'RocketPoolETH'
Set-TextValue $ws.Range("C50") 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextValue $ws.Range("D50") '2.022.43'
Set-TextValue $ws.Range("E50") '  -2.91%  '
Set-TextValue $ws.Range("D51") '940.71'
Set-TextValue $ws.Range("E51") '  -3.61%  '
